$d = $word.ActiveDocument

# 1) Megrendelő line: replace the name
$d.Content.Find.Execute("Megrendelő: Whastz the fuck you", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Megrendelő: 45sdfhgr5", 2)

# 2) Cím line: remove the trailing value, keep "Cím: "
$d.Content.Find.Execute("Cím: Xdddddd", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Cím: ", 2)

# 3) Date updates: all occurrences of 2024.07.01 -> 2024.07.02
$find = $d.Content.Find
$find.Execute("2024.07.01", $true, $false, $false, $false, $false,
               $true, 1, $false, "2024.07.02", 2)
while ($find.Execute("2024.07.01", $true, $false, $false, $false, $false,
                      $true, 1, $false, "2024.07.02", 2)) {
}
